$d = $word.ActiveDocument
$t = $d.Tables(1)

# "Sprint No." value cell: 1 -> 2
$sprintCell = $t.Cell(2, 4)
$sprintRange = $sprintCell.Range
$sprintRange.End = $sprintRange.End - 1
$sprintRange.Text = "2"

# "Review Date" value cell: 02/09/18 -> 02/21/18
$dateCell = $t.Cell(3, 2)
$dateRange = $dateCell.Range
$dateRange.End = $dateRange.End - 1
$dateRange.Text = "02/21/18"
